# Texas Holdem Bonus Simulator - add progress data
# Fills in the E-column (net winnings) figures for several hand rows that
# previously held placeholder zeros. The F (percentage) and G (scaled)
# columns are driven by existing shared formulas, so updating E is enough
# for them to recalculate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Winnings")

$values = @{
    "E93"  = -3281278936
    "E112" = 834290640
    "E113" = 10321042356
    "E133" = 9435360972
    "E134" = 18924544432
    "E156" = 19163347960
    "E157" = 28649104456
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Force recalculation so dependent formulas (F/G columns, I8, I11 totals)
# pick up the new values.
$excel.CalculateFull()

# Update the view to match where the user ended up scrolled/selected to.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 127
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E160").Select()
